# UN OCHA Briar updates
# Remove legacy ArcGIS fields (OBJECTID_1, OBJECTID__, Shape__, Shape_Length,
# Shape_Area, OBJECTID, Polygon) and clear the "<Null>" placeholder content
# left in empty fields (admin*RefName, admin*AltName1_en, admin*AltName2_en,
# validTo), then restore the intended column widths and a frozen header row
# on both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Admin0
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop the legacy ArcGIS columns: OBJECTID_1, OBJECTID__, Shape__ (A:C) and
# Shape_Length, Shape_Area, OBJECTID (originally L:N, now I:K after the
# first deletion shifts everything left).
$ws1.Range("A1:C1").EntireColumn.Delete()
$ws1.Range("I1:K1").EntireColumn.Delete()

# Clear the fields that only ever held the literal text "<Null>":
# admin0RefName, admin0AltName1_en, admin0AltName2_en, validTo
$ws1.Range("C2").Value = ""
$ws1.Range("D2").Value = ""
$ws1.Range("E2").Value = ""
$ws1.Range("H2").Value = ""

# Column widths (admin0Name_en, admin0Pcode, admin0RefName,
# admin0AltName1_en, admin0AltName2_en, date, validOn, validTo, Area_SqKm)
$ws1.Columns.Item(1).ColumnWidth = 24
$ws1.Columns.Item(2).ColumnWidth = 11
$ws1.Columns.Item(3).ColumnWidth = 13
$ws1.Columns.Item(4).ColumnWidth = 17
$ws1.Columns.Item(5).ColumnWidth = 17
$ws1.Columns.Item(6).ColumnWidth = 10
$ws1.Columns.Item(7).ColumnWidth = 10
$ws1.Columns.Item(8).ColumnWidth = 7
$ws1.Columns.Item(9).ColumnWidth = 15

# Freeze the header row.
$ws1.Activate()
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# Sheet 2: Admin1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Drop the legacy ArcGIS columns: OBJECTID_1, OBJECTID__, Shape__ (A:C) and
# Shape_Length, Shape_Area, OBJECTID (originally N:P, now K:M after the
# first deletion shifts everything left).
$ws2.Range("A1:C1").EntireColumn.Delete()
$ws2.Range("K1:M1").EntireColumn.Delete()

# Clear the fields that only ever held the literal text "<Null>":
# admin1RefName, admin1AltName1_en, admin1AltName2_en, validTo
$ws2.Range("C2").Value = ""
$ws2.Range("D2").Value = ""
$ws2.Range("E2").Value = ""
$ws2.Range("J2").Value = ""
$ws2.Range("C3").Value = ""
$ws2.Range("D3").Value = ""
$ws2.Range("E3").Value = ""
$ws2.Range("J3").Value = ""

# Column widths (admin1Name_en, admin1Pcode, admin1RefName,
# admin1AltName1_en, admin1AltName2_en, admin0Name_en, admin0Pcode, date,
# validOn, validTo, Area_SqKm)
$ws2.Columns.Item(1).ColumnWidth = 13
$ws2.Columns.Item(2).ColumnWidth = 11
$ws2.Columns.Item(3).ColumnWidth = 13
$ws2.Columns.Item(4).ColumnWidth = 17
$ws2.Columns.Item(5).ColumnWidth = 17
$ws2.Columns.Item(6).ColumnWidth = 24
$ws2.Columns.Item(7).ColumnWidth = 11
$ws2.Columns.Item(8).ColumnWidth = 10
$ws2.Columns.Item(9).ColumnWidth = 10
$ws2.Columns.Item(10).ColumnWidth = 7
$ws2.Columns.Item(11).ColumnWidth = 17

# Freeze the header row.
$ws2.Activate()
$ws2.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Leave the first sheet active/selected, matching the original workbook.
$ws1.Activate()
